$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column BF = 58. Rows 2..31 hold the game "Date" column, currently the
# mangled string "6-7-2012-13" (day-month-year2-year2 mashed together);
# the correct value is the real date "2013-06-07", stored as literal text
# (not a date serial) - same as the source data.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    if ($cell.Value2 -eq "6-7-2012-13") {
        # Assigning the literal string directly (or via Value2/Formula)
        # gets auto-parsed by Excel into a date serial number because the
        # text looks like a date. Route it through a formula that
        # evaluates to the literal text, then collapse the formula down
        # to its resulting value via copy / paste-values so the cell ends
        # up holding plain text "2013-06-07" with its original (default)
        # formatting untouched.
        $cell.Formula = '="2013-06-07"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = 0
